# Auto-generated edit script: Add data for 2022-08-21
# Updates 2022 year-to-date crime counts (column I, plus a couple H corrections)
# across the Citywide Totals, By Neighborhood summary, and individual neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 4566
$ws.Range("I3").Value = 4781
$ws.Range("H4").Value = 1672
$ws.Range("I4").Value = 1096
$ws.Range("I5").Value = 437
$ws.Range("I6").Value = 5205
$ws.Range("H7").Value = 25983
$ws.Range("I7").Value = 16085

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 129
$ws.Range("I6").Value = 109
$ws.Range("I7").Value = 512
$ws.Range("I8").Value = 977
$ws.Range("I11").Value = 243
$ws.Range("I12").Value = 38
$ws.Range("I15").Value = 185
$ws.Range("I17").Value = 22
$ws.Range("I19").Value = 453
$ws.Range("I23").Value = 150
$ws.Range("I29").Value = 1019
$ws.Range("I30").Value = 53
$ws.Range("I31").Value = 153
$ws.Range("I33").Value = 747
$ws.Range("I34").Value = 76
$ws.Range("I36").Value = 221
$ws.Range("I37").Value = 513
$ws.Range("I39").Value = 14
$ws.Range("I41").Value = 71
$ws.Range("I42").Value = 548
$ws.Range("I43").Value = 129
$ws.Range("I45").Value = 36
$ws.Range("I47").Value = 110
$ws.Range("I48").Value = 224
$ws.Range("I50").Value = 73
$ws.Range("I51").Value = 171
$ws.Range("I53").Value = 163
$ws.Range("I55").Value = 178
$ws.Range("I59").Value = 30
$ws.Range("H63").Value = 217
$ws.Range("I63").Value = 56
$ws.Range("I65").Value = 363
$ws.Range("I67").Value = 634
$ws.Range("I72").Value = 60
$ws.Range("I73").Value = 136
$ws.Range("I76").Value = 244
$ws.Range("I77").Value = 98
$ws.Range("I78").Value = 227
$ws.Range("I79").Value = 447
$ws.Range("I80").Value = 55
$ws.Range("I83").Value = 334
$ws.Range("I84").Value = 130
$ws.Range("I85").Value = 723
$ws.Range("I86").Value = 95
$ws.Range("I96").Value = 170
$ws.Range("I99").Value = 301
$ws.Range("H101").Value = 25983
$ws.Range("I101").Value = 16085

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I2").Value = 186
$ws.Range("I3").Value = 289
$ws.Range("I6").Value = 181
$ws.Range("I7").Value = 723

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 107
$ws.Range("I6").Value = 60
$ws.Range("I7").Value = 243

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 306
$ws.Range("I3").Value = 278
$ws.Range("I4").Value = 58
$ws.Range("I6").Value = 309
$ws.Range("I7").Value = 977

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I3").Value = 40
$ws.Range("I7").Value = 163

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I5").Value = 26
$ws.Range("I7").Value = 512

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I4").Value = 13
$ws.Range("I7").Value = 170

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("I6").Value = 16
$ws.Range("I7").Value = 53

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I6").Value = 138
$ws.Range("I7").Value = 513

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I4").Value = 22
$ws.Range("I6").Value = 81
$ws.Range("I7").Value = 301

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 221
$ws.Range("I7").Value = 634

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I3").Value = 43
$ws.Range("I7").Value = 153

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("I2").Value = 48
$ws.Range("I7").Value = 130

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I3").Value = 106
$ws.Range("I5").Value = 18
$ws.Range("I7").Value = 363

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 120
$ws.Range("I7").Value = 334

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I3").Value = 276
$ws.Range("I6").Value = 235
$ws.Range("I7").Value = 747

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 297
$ws.Range("I3").Value = 353
$ws.Range("I6").Value = 278
$ws.Range("I7").Value = 1019

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 168
$ws.Range("I3").Value = 133
$ws.Range("I6").Value = 123
$ws.Range("I7").Value = 453

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I3").Value = 43
$ws.Range("I7").Value = 224

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I2").Value = 52
$ws.Range("I7").Value = 244

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I2").Value = 48
$ws.Range("I6").Value = 24
$ws.Range("I7").Value = 109

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I3").Value = 25
$ws.Range("I7").Value = 71

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 142
$ws.Range("I3").Value = 187
$ws.Range("I4").Value = 43
$ws.Range("I6").Value = 157
$ws.Range("I7").Value = 548

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I2").Value = 52
$ws.Range("I6").Value = 86
$ws.Range("I7").Value = 227

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I3").Value = 59
$ws.Range("I6").Value = 55
$ws.Range("I7").Value = 178

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I2").Value = 42
$ws.Range("I3").Value = 53
$ws.Range("I7").Value = 150

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I2").Value = 132
$ws.Range("I3").Value = 144
$ws.Range("I7").Value = 447

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("I4").Value = 4
$ws.Range("I7").Value = 22

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I3").Value = 70
$ws.Range("I6").Value = 67
$ws.Range("I7").Value = 221

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("I6").Value = 15
$ws.Range("I7").Value = 76

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I3").Value = 35
$ws.Range("I4").Value = 10
$ws.Range("I7").Value = 110

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I6").Value = 65
$ws.Range("I7").Value = 185

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I4").Value = 15
$ws.Range("I7").Value = 73

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("I2").Value = 3
$ws.Range("I6").Value = 14

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I6").Value = 34
$ws.Range("I7").Value = 136

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("I2").Value = 14
$ws.Range("I7").Value = 30

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I2").Value = 42
$ws.Range("I7").Value = 129

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I4").Value = 47
$ws.Range("I7").Value = 95

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I5").Value = 3
$ws.Range("I7").Value = 171

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I3").Value = 22
$ws.Range("I7").Value = 129

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I6").Value = 34
$ws.Range("I7").Value = 60

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("I3").Value = 34
$ws.Range("I7").Value = 98

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("I3").Value = 8
$ws.Range("I7").Value = 36

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 55

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("I2").Value = 9
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 38
